# Running all the test cases
# Set Results column ("D") to "SKIP" for every test-case row on the
# "Test Cases" sheet (the header rows of each - possibly merged - block).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$rows = @(2,3,4,5,6,7,8,9,14,15,20,26,32,38,39,40,41,42,43,44,45)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = "SKIP"
}
